$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = "AAA BBB CCC"
$chars = $ws.Range("B10").Characters(5, 3)
$chars.Font.Bold = $true
$chars.Font.Color = 65535
